# Applies the edit described by the diff:
# 1. Rename sheet/tab "CopperA-HW40.xpc" to "CopperA"
# 2. Append a new data row (row 16) with averaged intensities for the
#    "HexGrid-60degTilt5degRes" scheme (A16=14, B16 label, C16:P16 values)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the worksheet / sheet tab
$ws.Name = "CopperA"

# 2a. Column A: index value, formatted like the cells above it (bold,
#     bordered, centered) -- copy formats only from A15 so no new style
#     definitions are introduced, then set the new value.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 14

# 2b. Column B: the scheme label (shared string already used in sheet)
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

# 2c. Columns C:P: the averaged intensity values
$values = @(
    0.9923753196759012,
    0.9884528854690178,
    0.995896977005419,
    0.9946912582886986,
    0.9923753196759012,
    0.9884528854690178,
    0.9921929612759364,
    0.9959692125204981,
    0.9946551405311592,
    0.9847566241747232,
    0.9923753196759012,
    0.9921749312372183,
    0.9928541101097592,
    0.9923737973676692
)

$col = 3
foreach ($v in $values) {
    $ws.Cells.Item(16, $col).Value = $v
    $col++
}
